$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the existing data, and a new blank column to
# its left so the original data block shifts from A1:O6 to B2:P7.
$ws.Rows.Item(1).Insert()
$ws.Columns.Item(1).Insert()

# Add the "Rat/Day" header label in the newly created A1 cell and make it bold.
$ws.Range("A1").Value = "Rat/Day"
$ws.Range("A1").Font.Bold = $true

# Leave the selection on the new header cell (matches the default/no-selection
# state Excel persists when the active cell is A1).
$ws.Range("A1").Select() | Out-Null
